$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.131.45"
$ws.Range("E2").Value = "  -4.54%  "
$ws.Range("D3").Value = "1.832.05"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'328.80"
$ws.Range("E5").Value = "  -2.77%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.4654"
$ws.Range("E7").Value = "  -2.19%  "
$ws.Range("D8").Value = "'0.3868"
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'46.16"
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07872"
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'0.9608"
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'22.00"
$ws.Range("E12").Value = "  -4.99%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.821.56"
$ws.Range("E13").Value = "  -4.38%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.670"
$ws.Range("E14").Value = "  -4.52%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.903"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.06856"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'1.0000"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "'86.63"
$ws.Range("E18").Value = "  -2.76%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000009958"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'16.67"
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.0000"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "28.143.13"
$ws.Range("E22").Value = "  -4.53%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.336"
$ws.Range("E23").Value = "  -3.11%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'11.03"
$ws.Range("E24").Value = "  -5.26%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.092"
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.035.29"
$ws.Range("E26").Value = "  -4.60%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'152.64"
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'19.21"
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'5.776"
$ws.Range("E29").Value = "  -10.94%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "'1.976"
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'117.20"
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.9366"
$ws.Range("E32").Value = "  -5.97%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.09260"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.309"
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.321"
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'3.349"
$ws.Range("E36").Value = "  -5.01%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05946"
$ws.Range("E37").Value = "  -7.41%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02148"
$ws.Range("E38").Value = "  -4.29%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.150"
$ws.Range("E39").Value = "  -4.03%  "
$ws.Range("D40").Value = "'7.687"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").Value = "'0.5601"
$ws.Range("E41").Value = "  -3.98%  "
$ws.Range("D42").Value = "'9.927"
$ws.Range("E42").Value = "  -5.75%  "
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").Value = "'1.228"
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("D45").Value = "'2.219"
$ws.Range("E45").Value = "  -8.38%  "
$ws.Range("E46").Value = "  -4.56%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.07050"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5271"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").Value = "'1.836"
$ws.Range("E49").Value = "  -6.03%  "
$ws.Range("D50").Value = "'111.75"
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("D51").Value = "'0.9995"
$ws.Range("E51").Value = "  -0.33%  "
